{"js": "// Applies the \"small edits to intro\" changes described by the diff:\n//  1. Reword the \"virtually any parameter...\" sentence (curse of dimensionality paragraph).\n//  2. \"Writing the package myself\" -> \"Writing the package\".\n//  3. Fix the GitHub URL casing and relocate the \"_GoBack\" bookmark to sit inside it\n//     (mirrors Word's own behaviour of leaving \"_GoBack\" at the last edited spot).\n\nasync function replaceFirst(body, searchText, newText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst doc = context.document;\n\n// 1) \"the possibilities are almost endless when attempting to fit models.\" ->\n//    \"the modeling possibilities are almost endless.\"\nawait replaceFirst(\n  body,\n  \"the possibilities are almost endless when attempting to fit models.\",\n  \"the modeling possibilities are almost endless.\"\n);\n\n// ... and \"the curse of dimensionality. For this reason,\" ->\n//     \"the curse of dimensionality while attempting to iteratively fit different parameter combinations to models. For this reason,\"\nawait replaceFirst(\n  body,\n  \"the curse of dimensionality. For this reason,\",\n  \"the curse of dimensionality while attempting to iteratively fit different parameter combinations to models. For this reason,\"\n);\n\n// 2) \"Writing the package myself has allowed for\" -> \"Writing the package has allowed for\"\nawait replaceFirst(\n  body,\n  \"Writing the package myself has allowed for\",\n  \"Writing the package has allowed for\"\n);\n\n// 3) Fix the GitHub link capitalization.\nawait replaceFirst(\n  body,\n  \"http://www.Github.com/Michael-Cowan/statmod.\",\n  \"http://www.github.com/michael-cowan/statmod.\"\n);\n\n// Move the hidden \"_GoBack\" bookmark so it sits where the link text was last edited\n// (right after \"michael-c\", before \"owan\"), matching Word's automatic behaviour.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst linkResults = body.search(\"michael-c\", { matchCase: false });\nlinkResults.load(\"items\");\nawait context.sync();\nif (linkResults.items.length > 0) {\n  const afterEdit = linkResults.items[0].getRange(Word.RangeLocation.end);\n  afterEdit.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Applies the \"small edits to intro\" changes described by the diff:\n#  1. Reword the \"virtually any parameter...\" sentence (curse of dimensionality paragraph).\n#  2. \"Writing the package myself\" -> \"Writing the package\".\n#  3. Fix the GitHub URL casing and relocate the \"_GoBack\" bookmark to sit inside it\n#     (mirrors Word's own behaviour of leaving \"_GoBack\" at the last edited spot).\n\n$d = $word.ActiveDocument\n\nfunction ReplaceOnce($findText, $replaceText) {\n    $rng = $d.Content\n    $null = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n\n# 1) \"the possibilities are almost endless when attempting to fit models.\" ->\n#    \"the modeling possibilities are almost endless.\"\nReplaceOnce \"the possibilities are almost endless when attempting to fit models.\" \"the modeling possibilities are almost endless.\"\n\n# ... and \"the curse of dimensionality. For this reason,\" ->\n#     \"the curse of dimensionality while attempting to iteratively fit different parameter combinations to models. For this reason,\"\nReplaceOnce \"the curse of dimensionality. For this reason,\" \"the curse of dimensionality while attempting to iteratively fit different parameter combinations to models. For this reason,\"\n\n# 2) \"Writing the package myself has allowed for\" -> \"Writing the package has allowed for\"\nReplaceOnce \"Writing the package myself has allowed for\" \"Writing the package has allowed for\"\n\n# 3) Fix the GitHub link capitalization.\nReplaceOnce \"http://www.Github.com/Michael-Cowan/statmod.\" \"http://www.github.com/michael-cowan/statmod.\"\n\n# Move the hidden \"_GoBack\" bookmark so it sits where the link text was last edited\n# (right after \"michael-c\", before \"owan\"), matching Word's automatic behaviour.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $oldBookmark = $d.Bookmarks.Item(\"_GoBack\")\n    $oldBookmark.Delete()\n}\n\n$editRng = $d.Content\n$found = $editRng.Find.Execute(\"michael-c\")\nif ($found) {\n    $editRng.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $editRng)\n}\n"}
